# Add a "Phone" column header and two new submission rows to the Card Data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("G1").Value = "Phone"

# Row 5 - entry without a phone number
$ws.Range("A5").Value = "JOHN DOE"
$ws.Range("B5").Value = "ahmed@gmail.com"
$ws.Range("C5").Value = "7007********5055"
$ws.Range("D5").Value = "12/∞"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "985"
$ws.Range("F5").Value = "11/13/2025, 12:03:28 PM"

# Row 6 - entry with a phone number
$ws.Range("A6").Value = "JOHN DOE"
$ws.Range("B6").Value = "johny@gmail.com"
$ws.Range("C6").Value = "7007********5055"
$ws.Range("D6").Value = "12/∞"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "985"
$ws.Range("F6").Value = "11/13/2025, 12:16:14 PM"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "+963993625082"
